$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5205.75
$ws.Range("I86").Value = 3850
$ws.Range("J86").Value = 5399.4287
$ws.Range("K86").Value = 3850
$ws.Range("L86").Value = 5399.4287
$ws.Range("M86").Value = -2727
$ws.Range("N86").Value = -7645.4287
$ws.Range("H89").Value = 5205.75
$ws.Range("I89").Value = 3850
$ws.Range("J89").Value = 5399.4287
$ws.Range("K89").Value = 19250
$ws.Range("L89").Value = 26997.1435
$ws.Range("M89").Value = -13634
$ws.Range("N89").Value = -38229.14350000001
$ws.Range("H98").Value = 1367.5555
$ws.Range("I98").Value = 1371.8572
$ws.Range("K98").Value = 1371.8572
$ws.Range("M98").Value = 126.1428000000001
$ws.Range("H100").Value = 653.7222
$ws.Range("I100").Value = 653.7222
$ws.Range("K100").Value = 653.7222
$ws.Range("M100").Value = -112.7222
$ws.Range("H101").Value = 766.3333
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H112").Value = 2425.4614
$ws.Range("J112").Value = 2556
$ws.Range("L112").Value = 7668
$ws.Range("N112").Value = -9884
$ws.Range("H118").Value = 9094.5
$ws.Range("I118").Value = 9094.5
$ws.Range("K118").Value = 27283.5
$ws.Range("M118").Value = -25626.5
$ws.Range("H122").Value = 1367.5555
$ws.Range("I122").Value = 1371.8572
$ws.Range("K122").Value = 4115.571599999999
$ws.Range("M122").Value = -1665.571599999999
$ws.Range("H137").Value = 8012.2383
$ws.Range("I137").Value = 2240.2
$ws.Range("J137").Value = 13259.546
$ws.Range("K137").Value = 6720.599999999999
$ws.Range("L137").Value = 39778.638
$ws.Range("M137").Value = -4170.599999999999
$ws.Range("N137").Value = -44878.638
$ws.Range("H138").Value = 5064.3296
$ws.Range("I138").Value = 6652.048
$ws.Range("J138").Value = 4588.014
$ws.Range("K138").Value = 19956.144
$ws.Range("L138").Value = 13764.042
$ws.Range("M138").Value = -14816.144
$ws.Range("N138").Value = -24044.042
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2007.1
$ws.Range("I2").Value = 2108.111
$ws.Range("J2").Value = 1098
$ws.Range("K2").Value = 2108.111
$ws.Range("L2").Value = 1098
$ws.Range("M2").Value = -1995.111
$ws.Range("N2").Value = -1324
$ws.Range("H5").Value = 162.75
$ws.Range("I5").Value = 150.33333
$ws.Range("K5").Value = 150.33333
$ws.Range("M5").Value = -38.33332999999999
$ws.Range("H32").Value = 15620.013
$ws.Range("I32").Value = 6213.104
$ws.Range("K32").Value = 6213.104
$ws.Range("M32").Value = -5926.104
$ws.Range("H35").Value = 1889.5
$ws.Range("I35").Value = 1889.5
$ws.Range("K35").Value = 1889.5
$ws.Range("M35").Value = -1483.5
$ws.Range("H116").Value = 2007.1
$ws.Range("I116").Value = 2108.111
$ws.Range("J116").Value = 1098
$ws.Range("K116").Value = 2108.111
$ws.Range("L116").Value = 1098
$ws.Range("M116").Value = 185.8890000000001
$ws.Range("N116").Value = -5686
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2007.1
$ws.Range("I3").Value = 2108.111
$ws.Range("J3").Value = 1098
$ws.Range("K3").Value = 2108.111
$ws.Range("L3").Value = 1098
$ws.Range("M3").Value = -1994.111
$ws.Range("N3").Value = -1326
$ws.Range("H4").Value = 162.75
$ws.Range("I4").Value = 150.33333
$ws.Range("K4").Value = 150.33333
$ws.Range("M4").Value = -35.33332999999999
$ws.Range("H134").Value = 4751.3335
$ws.Range("I134").Value = 2765
$ws.Range("K134").Value = 8295
$ws.Range("M134").Value = -5760
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 347.42856
$ws.Range("I16").Value = 356.4
$ws.Range("J16").Value = 325
$ws.Range("K16").Value = 356.4
$ws.Range("L16").Value = 325
$ws.Range("M16").Value = -69.39999999999998
$ws.Range("N16").Value = -899
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H113").Value = 347.42856
$ws.Range("I113").Value = 356.4
$ws.Range("J113").Value = 325
$ws.Range("K113").Value = 356.4
$ws.Range("L113").Value = 325
$ws.Range("M113").Value = 1813.6
$ws.Range("N113").Value = -4665
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2336.2727
$ws.Range("I34").Value = 1742.7142
$ws.Range("J34").Value = 3375
$ws.Range("K34").Value = 5228.142599999999
$ws.Range("L34").Value = 10125
$ws.Range("M34").Value = -5144.142599999999
$ws.Range("N34").Value = -10293
$ws.Range("H40").Value = 286.1111
$ws.Range("I40").Value = 273
$ws.Range("J40").Value = 312.33334
$ws.Range("K40").Value = 1092
$ws.Range("L40").Value = 1249.33336
$ws.Range("M40").Value = -1023
$ws.Range("N40").Value = -1387.33336
$ws.Range("H44").Value = 1154.2354
$ws.Range("I44").Value = 602.6667
$ws.Range("J44").Value = 1774.75
$ws.Range("K44").Value = 1808.0001
$ws.Range("L44").Value = 5324.25
$ws.Range("M44").Value = -1410.0001
$ws.Range("N44").Value = -6120.25
$ws.Range("H46").Value = 1667044.5
$ws.Range("I46").Value = 453.4
$ws.Range("K46").Value = 1360.2
$ws.Range("M46").Value = -1269.2
$ws.Range("H51").Value = 397
$ws.Range("I51").Value = 397
$ws.Range("K51").Value = 1191
$ws.Range("M51").Value = -731
$ws.Range("H131").Value = 1558.4762
$ws.Range("J131").Value = 1617.3158
$ws.Range("L131").Value = 4851.9474
$ws.Range("N131").Value = -14931.9474
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 958.75
$ws.Range("I2").Value = 274
$ws.Range("J2").Value = 2100
$ws.Range("K2").Value = 274
$ws.Range("L2").Value = 2100
$ws.Range("M2").Value = -161
$ws.Range("N2").Value = -2326
$ws.Range("H41").Value = 3408.6
$ws.Range("I41").Value = 3408.6
$ws.Range("K41").Value = 3408.6
$ws.Range("M41").Value = -3053.6
$ws.Range("H48").Value = 15000
$ws.Range("I48").Value = 15000
$ws.Range("K48").Value = 15000
$ws.Range("M48").Value = -14515
$ws.Range("H97").Value = 1765.32
$ws.Range("I97").Value = 1927.1578
$ws.Range("K97").Value = 1927.1578
$ws.Range("M97").Value = -1431.1578
$ws.Range("H107").Value = 94.5
$ws.Range("I107").Value = 94.5
$ws.Range("K107").Value = 94.5
$ws.Range("M107").Value = 1825.5
$ws.Range("H113").Value = 5223.5557
$ws.Range("I113").Value = 3999
$ws.Range("J113").Value = 5376.625
$ws.Range("K113").Value = 3999
$ws.Range("L113").Value = 5376.625
$ws.Range("M113").Value = -1829
$ws.Range("N113").Value = -9716.625
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -64900
$ws.Range("H126").Value = 4995.8335
$ws.Range("I126").Value = 4990
$ws.Range("K126").Value = 14970
$ws.Range("M126").Value = -12500
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4428.5713
$ws.Range("H61").Value = 3800.3076
$ws.Range("I61").Value = 3513.4783
$ws.Range("K61").Value = 3513.4783
$ws.Range("M61").Value = -3311.4783
$ws.Range("H100").Value = 2699.7778
$ws.Range("I100").Value = 2787.25
$ws.Range("K100").Value = 2787.25
$ws.Range("M100").Value = -2246.25
$ws.Range("H113").Value = 3800.3076
$ws.Range("I113").Value = 3513.4783
$ws.Range("K113").Value = 3513.4783
$ws.Range("M113").Value = -1343.4783
$ws.Range("H132").Value = 4668.6
$ws.Range("I132").Value = 3434.75
$ws.Range("J132").Value = 5807.5386
$ws.Range("K132").Value = 10304.25
$ws.Range("L132").Value = 17422.6158
$ws.Range("M132").Value = -7774.25
$ws.Range("N132").Value = -22482.6158
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2700
$ws.Range("I100").Value = 2640
$ws.Range("K100").Value = 5280
$ws.Range("M100").Value = -4739
$ws.Range("H113").Value = 1475.25
$ws.Range("I113").Value = 1475.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4425.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2255.75
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2259.4167
$ws.Range("I132").Value = 1719.5294
$ws.Range("J132").Value = 3570.5715
$ws.Range("K132").Value = 5158.5882
$ws.Range("L132").Value = 10711.7145
$ws.Range("M132").Value = -2628.5882
$ws.Range("N132").Value = -15771.7145
